$d = $word.ActiveDocument
$d.Content.Find.Execute("TABELA DE ORIENTAÇÃ ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "TABELA DE ORIENTAÇÃO", 2)
